# Weekly update for "Hortaliza, Terminal Hortofrutícola Agro Chillán - Zanahoria".
# A new weekly record is inserted at row 80 (pushing the existing historical
# rows 80-184 down to 81-185), and the new row 80 is populated with the
# latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 80; this shifts rows 80..184 down
# to 81..185 and grows the sheet's used range accordingly.
$ws.Rows(80).Insert()

# Populate the newly-inserted row 80 with the new weekly data point.
$ws.Cells.Item(80, 1).Value = 7
$ws.Cells.Item(80, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(80, 3).Value = "Ñuble"
$ws.Cells.Item(80, 4).Value = 44483
$ws.Cells.Item(80, 5).Value = 16
$ws.Cells.Item(80, 6).Value = 100114013
$ws.Cells.Item(80, 7).Value = "Zanahoria"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 160
$ws.Cells.Item(80, 11).Value = 7500
$ws.Cells.Item(80, 12).Value = 8000
$ws.Cells.Item(80, 13).Value = 7750
$ws.Cells.Item(80, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(80, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(80, 16).Value = 388
$ws.Cells.Item(80, 17).Value = 20
$ws.Cells.Item(80, 18).Value = "Hortaliza"
